# Updates profit-calculation figures on several sheets of the Leve profit
# tracker workbook (scheduled-runner refresh of market-board pricing data).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 43.846153
$ws.Range("I2").Value = 43.846153
$ws.Range("K2").Value = 43.846153
$ws.Range("M2").Value = 69.153847
# Row 9
$ws.Range("H9").Value = 1439.8572
$ws.Range("I9").Value = 1293.3334
$ws.Range("J9").Value = 1549.75
$ws.Range("K9").Value = 1293.3334
$ws.Range("L9").Value = 1549.75
$ws.Range("M9").Value = -1124.3334
$ws.Range("N9").Value = -1887.75
# Row 38
$ws.Range("H38").Value = 175
$ws.Range("I38").Value = 175
$ws.Range("K38").Value = 525
$ws.Range("M38").Value = -153
# Row 43
$ws.Range("H43").Value = 1373333
$ws.Range("J43").Value = 1373333
$ws.Range("L43").Value = 1373333
$ws.Range("N43").Value = -1373471
# Row 46
$ws.Range("H46").Value = 5000
$ws.Range("I46").Value = 5000
$ws.Range("K46").Value = 15000
$ws.Range("M46").Value = -14881
# Row 60
$ws.Range("H60").Value = 5000
$ws.Range("I60").Value = 5000
$ws.Range("K60").Value = 15000
$ws.Range("M60").Value = -14516
# Row 70
$ws.Range("H70").Value = 24307680
$ws.Range("J70").Value = 18520518
$ws.Range("L70").Value = 55561554
$ws.Range("N70").Value = -55562094
# Row 73
$ws.Range("H73").Value = 24307680
$ws.Range("J73").Value = 18520518
$ws.Range("L73").Value = 55561554
$ws.Range("N73").Value = -55563426
# Row 80
$ws.Range("H80").Value = 51289.9
$ws.Range("I80").Value = 50250
$ws.Range("J80").Value = 51549.875
$ws.Range("K80").Value = 150750
$ws.Range("L80").Value = 154649.625
$ws.Range("M80").Value = -149752
$ws.Range("N80").Value = -156645.625
# Row 83
$ws.Range("H83").Value = 51289.9
$ws.Range("I83").Value = 50250
$ws.Range("J83").Value = 51549.875
$ws.Range("K83").Value = 452250
$ws.Range("L83").Value = 463948.875
$ws.Range("M83").Value = -447258
$ws.Range("N83").Value = -473932.875
# Row 107
$ws.Range("H107").Value = 75001624
$ws.Range("I107").Value = 31252034
$ws.Range("K107").Value = 31252034
$ws.Range("M107").Value = -31250114
# Row 133
$ws.Range("H133").Value = 100762
$ws.Range("J133").Value = 100762
$ws.Range("L133").Value = 100762
$ws.Range("N133").Value = -110882
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# Row 141
$ws.Range("H141").Value = 1917
$ws.Range("I141").Value = 1952.0769
$ws.Range("K141").Value = 5856.2307
$ws.Range("M141").Value = -676.2307000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1605177
$ws.Range("I32").Value = 1763094.8
$ws.Range("K32").Value = 1763094.8
$ws.Range("M32").Value = -1762807.8
# Row 61
$ws.Range("H61").Value = 5530.8887
$ws.Range("J61").Value = 13218.889
$ws.Range("L61").Value = 13218.889
$ws.Range("N61").Value = -13642.889
# Row 74
$ws.Range("H74").Value = 115277.22
$ws.Range("I74").Value = 251874
$ws.Range("K74").Value = 251874
$ws.Range("M74").Value = -251000
# Row 77
$ws.Range("H77").Value = 115277.22
$ws.Range("I77").Value = 251874
$ws.Range("K77").Value = 1259370
$ws.Range("M77").Value = -1255002
# Row 136
$ws.Range("H136").Value = 5530.8887
$ws.Range("J136").Value = 13218.889
$ws.Range("L136").Value = 39656.667
$ws.Range("N136").Value = -44756.667

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 633.25
$ws.Range("J94").Value = 826.8570999999999
$ws.Range("L94").Value = 826.8570999999999
$ws.Range("N94").Value = -1728.8571
# Row 134
$ws.Range("H134").Value = 7942.241
$ws.Range("J134").Value = 10536.611
$ws.Range("L134").Value = 31609.833
$ws.Range("N134").Value = -36679.833

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6854.378
$ws.Range("I31").Value = 2673.84
$ws.Range("J31").Value = 12080.05
$ws.Range("K31").Value = 2673.84
$ws.Range("L31").Value = 12080.05
$ws.Range("M31").Value = -2378.84
$ws.Range("N31").Value = -12670.05
# Row 34
$ws.Range("H34").Value = 6854.378
$ws.Range("I34").Value = 2673.84
$ws.Range("J34").Value = 12080.05
$ws.Range("K34").Value = 2673.84
$ws.Range("L34").Value = 12080.05
$ws.Range("M34").Value = -2471.84
$ws.Range("N34").Value = -12484.05
# Row 55
$ws.Range("H55").Value = 29999
$ws.Range("I55").Value = 29999
$ws.Range("K55").Value = 29999
$ws.Range("M55").Value = -29684
# Row 107
$ws.Range("H107").Value = 1749.4667
$ws.Range("I107").Value = 1304
$ws.Range("K107").Value = 1304
$ws.Range("M107").Value = 616
# Row 132
$ws.Range("H132").Value = 5937.773
$ws.Range("I132").Value = 2603.6667
$ws.Range("J132").Value = 8246
$ws.Range("K132").Value = 7811.000100000001
$ws.Range("L132").Value = 24738
$ws.Range("M132").Value = -5281.000100000001
$ws.Range("N132").Value = -29798
# Row 134
$ws.Range("H134").Value = 7890.778
$ws.Range("J134").Value = 7944.0435
$ws.Range("L134").Value = 23832.1305
$ws.Range("N134").Value = -28902.1305

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 4651.75
$ws.Range("I34").Value = 114.666664
$ws.Range("J34").Value = 5698.769
$ws.Range("K34").Value = 343.999992
$ws.Range("L34").Value = 17096.307
$ws.Range("M34").Value = -259.999992
$ws.Range("N34").Value = -17264.307
# Row 39
$ws.Range("H39").Value = 10231.818
$ws.Range("I39").Value = 1125
$ws.Range("J39").Value = 15435.714
$ws.Range("K39").Value = 3375
$ws.Range("L39").Value = 46307.142
$ws.Range("M39").Value = -3081
$ws.Range("N39").Value = -46895.142
# Row 55
$ws.Range("H55").Value = 43341084
$ws.Range("I55").Value = 83334460
$ws.Range("J55").Value = 16678833
$ws.Range("K55").Value = 250003380
$ws.Range("L55").Value = 50036499
$ws.Range("M55").Value = -250003203
$ws.Range("N55").Value = -50036853

$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 72996.39999999999
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 72996.39999999999
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 72996.39999999999
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -74636.39999999999
# Row 97
$ws.Range("H97").Value = 1033.1569
$ws.Range("I97").Value = 953.46875
$ws.Range("J97").Value = 1167.3684
$ws.Range("K97").Value = 953.46875
$ws.Range("L97").Value = 1167.3684
$ws.Range("M97").Value = -457.46875
$ws.Range("N97").Value = -2159.3684
# Row 126
$ws.Range("H126").Value = 31252174
$ws.Range("I126").Value = 62501910
$ws.Range("K126").Value = 187505730
$ws.Range("M126").Value = -187503260

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 3646.0256
$ws.Range("I122").Value = 2823.818
$ws.Range("K122").Value = 8471.454000000002
$ws.Range("M122").Value = -6021.454000000002
# Row 135
$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 714.26666
$ws.Range("I107").Value = 517.9167
$ws.Range("K107").Value = 1553.7501
$ws.Range("M107").Value = 366.2499
# Row 132
$ws.Range("H132").Value = 71480000
$ws.Range("J132").Value = 84996.336
$ws.Range("L132").Value = 254989.008
$ws.Range("N132").Value = -260049.008
# Row 139
$ws.Range("H139").Value = 88643.57000000001
$ws.Range("J139").Value = 88643.57000000001
$ws.Range("L139").Value = 88643.57000000001
$ws.Range("N139").Value = -98923.57000000001
